$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 5.830900000000002
$ws.Range("B14").Value = 6.375300000000006
$ws.Range("B21").Value = 9.509199999999995
$ws.Range("C22").Value = -11.6915
$ws.Range("B23").Value = 8.969199999999997
$ws.Range("C24").Value = -13.7465
$ws.Range("B25").Value = 5.622000000000002
$ws.Range("B26").Value = 4.968300000000003
$ws.Range("C28").Value = -14.0269
$ws.Range("B29").Value = 5.368500000000004
$ws.Range("C36").Value = -11.9482
$ws.Range("C45").Value = -13.83859999999999
$ws.Range("C48").Value = -12.2521
$ws.Range("C49").Value = -13.76299999999999
$ws.Range("C52").Value = -10.73409999999999
$ws.Range("B53").Value = 5.691400000000001
$ws.Range("C53").Value = -10.86220000000001
$ws.Range("C54").Value = -13.65859999999999
$ws.Range("B57").Value = 4.572099999999994
$ws.Range("B59").Value = 5.007
$ws.Range("B69").Value = 5.661999999999993
$ws.Range("C70").Value = -12.5156
$ws.Range("B79").Value = 9.399900000000004
$ws.Range("B83").Value = 5.512799999999996
$ws.Range("C86").Value = -13.5531
$ws.Range("C87").Value = -12.8106
$ws.Range("C89").Value = -13.2964
$ws.Range("B91").Value = 5.907500000000002
$ws.Range("B93").Value = 5.495299999999997
$ws.Range("C101").Value = -13.4387
$ws.Range("B103").Value = 5.613100000000006
